$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the value in cell D1 (new classes for accept env)
$ws.Range("D1").Value = 511377915
